$wb = $excel.ActiveWorkbook

# --- Sheet2 (Column Name / Value reference table) ---------------------------
# Populate new rows 2-5 FIRST (in top-to-bottom, left-to-right reading order)
# so that newly introduced shared strings are appended to the shared-string
# table in the same order the final workbook expects.
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = "Quote Order Number"
$ws2.Range("B2").Value = "2021-78 — ETCC UFT licence quote"
$ws2.Range("D2").Value = "shipping Number"
$ws2.Range("E2").Value = "2021-45 — ETCC UFT licence Shipping"

$ws2.Range("A3").Value = "Quote Status"
$ws2.Range("B3").Value = "DRAFT"
$ws2.Range("D3").Value = "shipping Status"
$ws2.Range("E3").Value = "SHIPPED"

$ws2.Range("A4").Value = "Sales Order Number"
$ws2.Range("B4").Value = "2021-53 — ETCC UFT licence Sales Order"
$ws2.Range("D4").Value = "Invoice Number"
$ws2.Range("E4").Value = "2021-53 — ETCC UFT licence Sales Order"

$ws2.Range("A5").Value = "SO Status"
$ws2.Range("B5").Value = "Ordered"
$ws2.Range("D5").Value = "Invoice Status"
$ws2.Range("E5").Value = "Partially Shipped & Invoiced"

# --- Sheet1 (Test Case data) -------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("M1").Value = "Order Number"

$ws1.Range("A2").Value = "TC001_Creating_Qoute_Converting_To_SalesOrder"
$ws1.Range("B2").Value = "SalesOrderTest"
$ws1.Range("C2").Value = "iDeliver@1"
$ws1.Range("D2").Value = "ETCC UFT licence"
$ws1.Range("E2").Value = "5 UFT licence for ETCC"
$ws1.Range("F2").Value = "etcc"

$ws1.Range("A3").Value = "TC002_Converting_SalesOrder_To_Invoicing"
$ws1.Range("B3").Value = "ShipmentTestUser"
$ws1.Range("C3").Value = "iDeliver@1"
$ws1.Range("D3").Value = "ETCC UFT licence"
$ws1.Range("E3").Value = "5 UFT licence for ETCC"
$ws1.Range("F3").Value = "24/7 Couriers"
$ws1.Range("M3").Value = "TC001_Creating_Qoute_Converting_To_SalesOrder;Sales Order Number"

$ws1.Range("A4").Value = "TC003_Invoice_To_Payment"
$ws1.Range("B4").Value = "ProcessBillingUser"
$ws1.Range("C4").Value = "iDeliver@1"
$ws1.Range("D4").Value = "ETCC UFT licence"
$ws1.Range("E4").Value = "5 UFT licence for ETCC"
$ws1.Range("F4").Value = "24/7 Couriers"
$ws1.Range("M4").Value = "TC002_Converting_SalesOrder_To_Invoicing;Invoice Number"

# --- Selection state (both sheets end up with E4 selected) ------------------
$ws1.Activate()
$ws1.Range("E4").Select()

$ws2.Activate()
$ws2.Range("E4").Select()

$ws1.Activate()
